$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$tbl = $ws.ListObjects.Item(1)
$tbl.ListColumns.Item("Buy Value").Name = "Buy Value in GBP"
$tbl.ListColumns.Item("Sell Value").Name = "Sell Value in GBP"
$tbl.ListColumns.Item("Fee Value").Name = "Fee Value in GBP"
